$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$data = @(
    @("Payton Pritchard", "PG,SG", "Boston Celtics"),
    @("Ty Jerome", "PG,SG", "Cleveland Cavaliers"),
    @("Jalen Green", "PG,SG", "Houston Rockets"),
    @("Chris Paul", "PG", "San Antonio Spurs"),
    @("Jaylen Brown", "SG,SF", "Boston Celtics"),
    @("Brandon Clarke", "PF,C", "Memphis Grizzlies"),
    @("Ayo Dosunmu", "PG,SG,SF", "Chicago Bulls"),
    @("Pascal Siakam", "SF,PF,C", "Indiana Pacers"),
    @("Nikola Jokic", "C", "Denver Nuggets"),
    @("Rudy Gobert", "C", "Minnesota Timberwolves"),
    @("Jakob Poeltl", "C", "Toronto Raptors"),
    @("Stephon Castle", "PG,SG", "San Antonio Spurs"),
    @("Dillon Brooks", "SG,SF", "Houston Rockets"),
    @("Paolo Banchero", "SF,PF", "Orlando Magic"),
    @("Deni Avdija", "SF,PF", "Portland Trail Blazers"),
    @("Chet Holmgren", "PF,C", "Oklahoma City Thunder"),
    @("Jalen Suggs", "PG,SG", "Orlando Magic"),
    @("Russell Westbrook", "PG,SG", "Denver Nuggets")
)

for ($i = 0; $i -lt $data.Length; $i++) {
    $row = $i + 2
    $ws.Cells.Item($row, 1).Value = $data[$i][0]
}
for ($i = 0; $i -lt $data.Length; $i++) {
    $row = $i + 2
    $ws.Cells.Item($row, 2).Value = $data[$i][1]
}
for ($i = 0; $i -lt $data.Length; $i++) {
    $row = $i + 2
    $ws.Cells.Item($row, 3).Value = $data[$i][2]
}
